$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 41
$ws.Range("H41").Value = 665.4583
$ws.Range("I41").Value = 587.6
$ws.Range("J41").Value = 721.0714
$ws.Range("K41").Value = 587.6
$ws.Range("L41").Value = 721.0714
$ws.Range("M41").Value = -147.6
$ws.Range("N41").Value = -1601.0714
# row 62
$ws.Range("H62").Value = 2753.2222
$ws.Range("I62").Value = 2000
$ws.Range("K62").Value = 2000
$ws.Range("M62").Value = -1376
# row 65
$ws.Range("H65").Value = 2753.2222
$ws.Range("I65").Value = 2000
$ws.Range("K65").Value = 10000
$ws.Range("M65").Value = -6880
# row 86
$ws.Range("H86").Value = 54811.05
$ws.Range("I86").Value = 85507.75
$ws.Range("K86").Value = 85507.75
$ws.Range("M86").Value = -84384.75
# row 89
$ws.Range("H89").Value = 54811.05
$ws.Range("I89").Value = 85507.75
$ws.Range("K89").Value = 427538.75
$ws.Range("M89").Value = -421922.75
# row 106
$ws.Range("H106").Value = 1800
$ws.Range("I106").Value = 1800
$ws.Range("K106").Value = 1800
$ws.Range("M106").Value = -1169
# row 129
$ws.Range("H129").Value = 1961.2959
$ws.Range("I129").Value = 6047.3887
$ws.Range("J129").Value = 1041.925
$ws.Range("K129").Value = 18142.1661
$ws.Range("L129").Value = 3125.775
$ws.Range("M129").Value = -13142.1661
$ws.Range("N129").Value = -13125.775
# row 141
$ws.Range("H141").Value = 2632.2307
$ws.Range("I141").Value = 2538.0908
$ws.Range("J141").Value = 3150
$ws.Range("K141").Value = 7614.2724
$ws.Range("L141").Value = 9450
$ws.Range("M141").Value = -2434.2724
$ws.Range("N141").Value = -19810

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 30810.262
$ws.Range("I32").Value = 6782.9653
$ws.Range("J32").Value = 495338
$ws.Range("K32").Value = 6782.9653
$ws.Range("L32").Value = 495338
$ws.Range("M32").Value = -6495.9653
$ws.Range("N32").Value = -495912
# row 45
$ws.Range("H45").Value = 101539.8
$ws.Range("I45").Value = 126560.875
$ws.Range("K45").Value = 126560.875
$ws.Range("M45").Value = -126183.875

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 35
$ws.Range("H35").Value = 0
$ws.Range("I35").Value = 0
$ws.Range("K35").Value = 0
$ws.Range("M35").ClearContents()
# row 132
$ws.Range("H132").Value = 3979.8235
$ws.Range("I132").Value = 3600.0435
$ws.Range("J132").Value = 4773.909
$ws.Range("K132").Value = 10800.1305
$ws.Range("L132").Value = 14321.727
$ws.Range("M132").Value = -8270.130500000001
$ws.Range("N132").Value = -19381.727

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 16.875
$ws.Range("I2").Value = 18
$ws.Range("J2").Value = 15.75
$ws.Range("K2").Value = 108
$ws.Range("L2").Value = 94.5
$ws.Range("M2").Value = 5
$ws.Range("N2").Value = -320.5
# row 5
$ws.Range("H5").Value = 1255.5227
$ws.Range("I5").Value = 1153.7693
$ws.Range("K5").Value = 3461.3079
$ws.Range("M5").Value = -3349.3079
# row 23
$ws.Range("H23").Value = 997.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 997.5
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 2992.5
$ws.Range("N23").Value = -3462.5
$ws.Range("M23").ClearContents()
# row 86
$ws.Range("H86").Value = 816
$ws.Range("I86").Value = 480
$ws.Range("J86").Value = 900
$ws.Range("K86").Value = 1440
$ws.Range("L86").Value = 2700
$ws.Range("M86").Value = -254
$ws.Range("N86").Value = -5072
# row 89
$ws.Range("H89").Value = 816
$ws.Range("I89").Value = 480
$ws.Range("J89").Value = 900
$ws.Range("K89").Value = 4320
$ws.Range("L89").Value = 8100
$ws.Range("M89").Value = 1608
$ws.Range("N89").Value = -19956
# row 97
$ws.Range("H97").Value = 1672
$ws.Range("I97").Value = 1360
$ws.Range("J97").Value = 2452
$ws.Range("K97").Value = 4080
$ws.Range("L97").Value = 7356
$ws.Range("M97").Value = -3584
$ws.Range("N97").Value = -8348
# row 132
$ws.Range("H132").Value = 1659.6666
$ws.Range("I132").Value = 799.8889
$ws.Range("J132").Value = 1946.2593
$ws.Range("K132").Value = 7199.0001
$ws.Range("L132").Value = 17516.3337
$ws.Range("M132").Value = -4669.0001
$ws.Range("N132").Value = -22576.3337
# row 135
$ws.Range("H135").Value = 1255.5227
$ws.Range("I135").Value = 1153.7693
$ws.Range("K135").Value = 10383.9237
$ws.Range("M135").Value = -7848.923699999999

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 102
$ws.Range("H102").Value = 3329.375
$ws.Range("I102").Value = 2480
$ws.Range("J102").Value = 3839
$ws.Range("K102").Value = 2480
$ws.Range("L102").Value = 3839
$ws.Range("M102").Value = -858
$ws.Range("N102").Value = -7083
# row 122
$ws.Range("H122").Value = 796.931
$ws.Range("I122").Value = 605.2857
$ws.Range("K122").Value = 1815.8571
$ws.Range("M122").Value = 634.1428999999998

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 2535.1428
$ws.Range("I7").Value = 2065.6667
$ws.Range("J7").Value = 3161.111
$ws.Range("K7").Value = 2065.6667
$ws.Range("L7").Value = 3161.111
$ws.Range("M7").Value = -1953.6667
$ws.Range("N7").Value = -3385.111
# row 82
$ws.Range("H82").Value = 1760.7894
$ws.Range("I82").Value = 1409.2727
$ws.Range("J82").Value = 2244.125
$ws.Range("K82").Value = 1409.2727
$ws.Range("L82").Value = 2244.125
$ws.Range("M82").Value = -1048.2727
$ws.Range("N82").Value = -2966.125
# row 85
$ws.Range("H85").Value = 1760.7894
$ws.Range("I85").Value = 1409.2727
$ws.Range("J85").Value = 2244.125
$ws.Range("K85").Value = 1409.2727
$ws.Range("L85").Value = 2244.125
$ws.Range("M85").Value = -161.2727
$ws.Range("N85").Value = -4740.125
# row 126
$ws.Range("H126").Value = 2535.1428
$ws.Range("I126").Value = 2065.6667
$ws.Range("J126").Value = 3161.111
$ws.Range("K126").Value = 6197.000100000001
$ws.Range("L126").Value = 9483.332999999999
$ws.Range("M126").Value = -3727.000100000001
$ws.Range("N126").Value = -14423.333

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 47
$ws.Range("H47").Value = 32556.334
$ws.Range("J47").Value = 32556.334
$ws.Range("L47").Value = 32556.334
$ws.Range("N47").Value = -33700.334
# row 107
$ws.Range("H107").Value = 100431.2
$ws.Range("I107").Value = 288.66666
$ws.Range("J107").Value = 250645
$ws.Range("K107").Value = 865.9999799999999
$ws.Range("L107").Value = 751935
$ws.Range("M107").Value = 1054.00002
$ws.Range("N107").Value = -755775
# row 113
$ws.Range("H113").Value = 822.05554
$ws.Range("I113").Value = 633.3333
$ws.Range("J113").Value = 1010.7778
$ws.Range("K113").Value = 1899.9999
$ws.Range("L113").Value = 3032.3334
$ws.Range("M113").Value = 270.0001
$ws.Range("N113").Value = -7372.3334
